# Edit script generated to apply diff changes to germany_2-bundesliga_2023-2024 workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Append two new match rows (128, 129) at the end of the table.
#    First clone formatting/styles from the last existing row (127),
#    then overwrite every cell with the correct values.
# ------------------------------------------------------------------
$ws.Range("A127:V127").Copy($ws.Range("A128:V128"))
$ws.Range("A127:V127").Copy($ws.Range("A129:V129"))

# Row 128
$ws.Range("A128").Value = 127
$ws.Range("B128").Value = 'germany'
$ws.Range("C128").Value = '2-bundesliga'
$ws.Range("D128").Value = '2023-2024'
$ws.Range("E128").Value = 45261.77083333334
$ws.Range("F128").Value = 'Schalke'
$ws.Range("G128").Value = 4
$ws.Range("H128").Value = 'VfL Osnabruck'
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 1.69
$ws.Range("K128").Value = '25/11/2023 20:42'
$ws.Range("L128").Value = 1.53
$ws.Range("M128").Value = '01/12/2023 18:29'
$ws.Range("N128").Value = 4.57
$ws.Range("O128").Value = '25/11/2023 20:42'
$ws.Range("P128").Value = 4.51
$ws.Range("Q128").Value = '01/12/2023 18:29'
$ws.Range("R128").Value = 4.27
$ws.Range("S128").Value = '25/11/2023 20:42'
$ws.Range("T128").Value = 6.32
$ws.Range("U128").Value = '01/12/2023 18:28'
$ws.Range("V128").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/schalke-vfl-osnabruck/vVBy4buB/'

# Row 129
$ws.Range("A129").Value = 128
$ws.Range("B129").Value = 'germany'
$ws.Range("C129").Value = '2-bundesliga'
$ws.Range("D129").Value = '2023-2024'
$ws.Range("E129").Value = 45261.77083333334
$ws.Range("F129").Value = 'St. Pauli'
$ws.Range("G129").Value = 2
$ws.Range("H129").Value = 'Hamburger SV'
$ws.Range("I129").Value = 2
$ws.Range("J129").Value = 2.02
$ws.Range("K129").Value = '25/11/2023 13:13'
$ws.Range("L129").Value = 1.87
$ws.Range("M129").Value = '01/12/2023 18:25'
$ws.Range("N129").Value = 3.96
$ws.Range("O129").Value = '25/11/2023 13:13'
$ws.Range("P129").Value = 3.88
$ws.Range("Q129").Value = '01/12/2023 18:25'
$ws.Range("R129").Value = 3.42
$ws.Range("S129").Value = '25/11/2023 13:13'
$ws.Range("T129").Value = 4.21
$ws.Range("U129").Value = '01/12/2023 18:25'
$ws.Range("V129").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/st-pauli-hamburger/YT4p2xAN/'

# ------------------------------------------------------------------
# 2) Re-shuffle the match details (home/away teams, scores, odds, url)
#    for the rows whose fixtures were re-ordered between scrapes.
#    The row Index (A) and match date/time (E) stay put; only the
#    F:V block with the actual fixture content moves between rows.
# ------------------------------------------------------------------
# Row 3
$ws.Range("F3").Value = 'Wehen'
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 'Magdeburg'
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3.22
$ws.Range("K3").Value = '08/07/2023 18:28'
$ws.Range("L3").Value = 3.09
$ws.Range("M3").Value = '29/07/2023 12:33'
$ws.Range("N3").Value = 3.69
$ws.Range("O3").Value = '08/07/2023 18:28'
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Value = '29/07/2023 12:58'
$ws.Range("R3").Value = 2.18
$ws.Range("S3").Value = '08/07/2023 18:28'
$ws.Range("T3").Value = 2.28
$ws.Range("U3").Value = '29/07/2023 12:33'
$ws.Range("V3").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/wehen-magdeburg/x0US8naJ/'

# Row 4
$ws.Range("F4").Value = 'Hannover'
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 'Elversberg'
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 1.98
$ws.Range("K4").Value = '08/07/2023 18:28'
$ws.Range("L4").Value = 1.96
$ws.Range("M4").Value = '29/07/2023 12:58'
$ws.Range("N4").Value = 3.78
$ws.Range("O4").Value = '08/07/2023 18:28'
$ws.Range("P4").Value = 4.09
$ws.Range("Q4").Value = '29/07/2023 12:58'
$ws.Range("R4").Value = 3.7
$ws.Range("S4").Value = '08/07/2023 18:28'
$ws.Range("T4").Value = 3.65
$ws.Range("U4").Value = '29/07/2023 12:58'
$ws.Range("V4").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hannover-elversberg/pWdGDApl/'

# Row 5
$ws.Range("F5").Value = 'Kaiserslautern'
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 'St. Pauli'
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 3.1
$ws.Range("K5").Value = '08/07/2023 18:28'
$ws.Range("L5").Value = 3.27
$ws.Range("M5").Value = '29/07/2023 12:59'
$ws.Range("N5").Value = 3.57
$ws.Range("O5").Value = '08/07/2023 18:28'
$ws.Range("P5").Value = 3.12
$ws.Range("Q5").Value = '29/07/2023 12:57'
$ws.Range("R5").Value = 2.33
$ws.Range("S5").Value = '08/07/2023 18:28'
$ws.Range("T5").Value = 2.5
$ws.Range("U5").Value = '29/07/2023 12:59'
$ws.Range("V5").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/kaiserslautern-st-pauli/j9cCEjVs/'

# Row 6
$ws.Range("F6").Value = 'VfL Osnabruck'
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 'Karlsruher SC'
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 2.75
$ws.Range("K6").Value = '08/07/2023 18:28'
$ws.Range("L6").Value = 3.19
$ws.Range("M6").Value = '29/07/2023 12:59'
$ws.Range("N6").Value = 3.58
$ws.Range("O6").Value = '08/07/2023 18:28'
$ws.Range("P6").Value = 3.84
$ws.Range("Q6").Value = '29/07/2023 12:59'
$ws.Range("R6").Value = 2.52
$ws.Range("S6").Value = '08/07/2023 18:28'
$ws.Range("T6").Value = 2.2
$ws.Range("U6").Value = '29/07/2023 12:59'
$ws.Range("V6").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/vfl-osnabruck-karlsruher/OdQO9SqD/'

# Row 8
$ws.Range("F8").Value = 'Hansa Rostock'
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 'Nurnberg'
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2.42
$ws.Range("K8").Value = '08/07/2023 18:28'
$ws.Range("L8").Value = 2.35
$ws.Range("M8").Value = '30/07/2023 13:28'
$ws.Range("N8").Value = 3.33
$ws.Range("O8").Value = '08/07/2023 18:28'
$ws.Range("P8").Value = 3.25
$ws.Range("Q8").Value = '30/07/2023 13:27'
$ws.Range("R8").Value = 3.07
$ws.Range("S8").Value = '08/07/2023 18:28'
$ws.Range("T8").Value = 3.39
$ws.Range("U8").Value = '30/07/2023 13:28'
$ws.Range("V8").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hansa-rostock-nurnberg/jqSGBlF0/'

# Row 9
$ws.Range("F9").Value = 'Greuther Furth'
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 'Paderborn'
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 2.53
$ws.Range("K9").Value = '08/07/2023 18:28'
$ws.Range("L9").Value = 2.57
$ws.Range("M9").Value = '30/07/2023 13:26'
$ws.Range("N9").Value = 3.61
$ws.Range("O9").Value = '08/07/2023 18:28'
$ws.Range("P9").Value = 3.7
$ws.Range("Q9").Value = '30/07/2023 13:26'
$ws.Range("R9").Value = 2.78
$ws.Range("S9").Value = '08/07/2023 18:28'
$ws.Range("T9").Value = 2.71
$ws.Range("U9").Value = '30/07/2023 13:28'
$ws.Range("V9").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/greuther-furth-paderborn/OMeKCUaf/'

# Row 10
$ws.Range("F10").Value = 'Braunschweig'
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 'Holstein Kiel'
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 2.5
$ws.Range("K10").Value = '08/07/2023 18:28'
$ws.Range("L10").Value = 2.35
$ws.Range("M10").Value = '30/07/2023 13:29'
$ws.Range("N10").Value = 3.8
$ws.Range("O10").Value = '08/07/2023 18:28'
$ws.Range("P10").Value = 3.72
$ws.Range("Q10").Value = '30/07/2023 13:28'
$ws.Range("R10").Value = 2.67
$ws.Range("S10").Value = '08/07/2023 18:28'
$ws.Range("T10").Value = 2.99
$ws.Range("U10").Value = '30/07/2023 13:29'
$ws.Range("V10").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/braunschweig-holstein-kiel/UZQKA8U6/'

# Row 17
$ws.Range("F17").Value = 'Karlsruher SC'
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = 'Hamburger SV'
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 2.93
$ws.Range("K17").Value = '30/07/2023 13:42'
$ws.Range("L17").Value = 2.48
$ws.Range("M17").Value = '06/08/2023 13:27'
$ws.Range("N17").Value = 3.97
$ws.Range("O17").Value = '30/07/2023 13:42'
$ws.Range("P17").Value = 4.14
$ws.Range("Q17").Value = '06/08/2023 13:28'
$ws.Range("R17").Value = 2.28
$ws.Range("S17").Value = '30/07/2023 13:42'
$ws.Range("T17").Value = 2.62
$ws.Range("U17").Value = '06/08/2023 13:28'
$ws.Range("V17").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/karlsruher-hamburger/O0G2aqyb/'

# Row 18
$ws.Range("F18").Value = 'Nurnberg'
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 'Hannover'
$ws.Range("I18").Value = 2
$ws.Range("J18").Value = 2
$ws.Range("K18").Value = '30/07/2023 13:42'
$ws.Range("L18").Value = 2.48
$ws.Range("M18").Value = '06/08/2023 13:26'
$ws.Range("N18").Value = 3.78
$ws.Range("O18").Value = '30/07/2023 13:42'
$ws.Range("P18").Value = 3.45
$ws.Range("Q18").Value = '06/08/2023 13:28'
$ws.Range("R18").Value = 3.73
$ws.Range("S18").Value = '30/07/2023 13:42'
$ws.Range("T18").Value = 2.98
$ws.Range("U18").Value = '06/08/2023 13:29'
$ws.Range("V18").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/nurnberg-hannover/bJ8FdsLG/'

# Row 19
$ws.Range("F19").Value = 'Magdeburg'
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 'Braunschweig'
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 2.05
$ws.Range("K19").Value = '30/07/2023 13:42'
$ws.Range("L19").Value = 1.78
$ws.Range("M19").Value = '06/08/2023 13:23'
$ws.Range("N19").Value = 3.87
$ws.Range("O19").Value = '30/07/2023 13:42'
$ws.Range("P19").Value = 4.02
$ws.Range("Q19").Value = '06/08/2023 13:27'
$ws.Range("R19").Value = 3.49
$ws.Range("S19").Value = '30/07/2023 13:42'
$ws.Range("T19").Value = 4.53
$ws.Range("U19").Value = '06/08/2023 13:29'
$ws.Range("V19").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/magdeburg-braunschweig/KS9BcN5A/'

# Row 29
$ws.Range("F29").Value = 'Schalke'
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 'Holstein Kiel'
$ws.Range("I29").Value = 2
$ws.Range("J29").Value = 1.65
$ws.Range("K29").Value = '20/08/2023 15:12'
$ws.Range("L29").Value = 2.02
$ws.Range("M29").Value = '25/08/2023 18:29'
$ws.Range("N29").Value = 4.41
$ws.Range("O29").Value = '20/08/2023 15:12'
$ws.Range("P29").Value = 4.09
$ws.Range("Q29").Value = '25/08/2023 18:29'
$ws.Range("R29").Value = 4.78
$ws.Range("S29").Value = '20/08/2023 15:12'
$ws.Range("T29").Value = 3.45
$ws.Range("U29").Value = '25/08/2023 18:29'
$ws.Range("V29").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/schalke-holstein-kiel/2BR18vKS/'

# Row 30
$ws.Range("F30").Value = 'Paderborn'
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 'Kaiserslautern'
$ws.Range("I30").Value = 2
$ws.Range("J30").Value = 1.8
$ws.Range("K30").Value = '19/08/2023 12:12'
$ws.Range("L30").Value = 2.17
$ws.Range("M30").Value = '25/08/2023 18:25'
$ws.Range("N30").Value = 4.19
$ws.Range("O30").Value = '19/08/2023 12:12'
$ws.Range("P30").Value = 3.74
$ws.Range("Q30").Value = '25/08/2023 18:24'
$ws.Range("R30").Value = 4.04
$ws.Range("S30").Value = '19/08/2023 12:12'
$ws.Range("T30").Value = 3.33
$ws.Range("U30").Value = '25/08/2023 18:28'
$ws.Range("V30").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/paderborn-kaiserslautern/jyCy2xcd/'

# Row 44
$ws.Range("F44").Value = 'VfL Osnabruck'
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 'Elversberg'
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = 2.35
$ws.Range("K44").Value = '27/08/2023 17:12'
$ws.Range("L44").Value = 2.32
$ws.Range("M44").Value = '03/09/2023 13:02'
$ws.Range("N44").Value = 3.78
$ws.Range("O44").Value = '27/08/2023 17:12'
$ws.Range("P44").Value = 4.01
$ws.Range("Q44").Value = '03/09/2023 13:04'
$ws.Range("R44").Value = 2.93
$ws.Range("S44").Value = '27/08/2023 17:12'
$ws.Range("T44").Value = 2.88
$ws.Range("U44").Value = '03/09/2023 13:04'
$ws.Range("V44").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/vfl-osnabruck-elversberg/djbBOjAQ/'

# Row 46
$ws.Range("F46").Value = 'Greuther Furth'
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 'Hannover'
$ws.Range("I46").Value = 3
$ws.Range("J46").Value = 1.72
$ws.Range("K46").Value = '27/08/2023 12:42'
$ws.Range("L46").Value = 2.03
$ws.Range("M46").Value = '03/09/2023 13:27'
$ws.Range("N46").Value = 4.27
$ws.Range("O46").Value = '27/08/2023 12:42'
$ws.Range("P46").Value = 3.76
$ws.Range("Q46").Value = '03/09/2023 13:18'
$ws.Range("R46").Value = 4.58
$ws.Range("S46").Value = '27/08/2023 12:42'
$ws.Range("T46").Value = 3.7
$ws.Range("U46").Value = '03/09/2023 13:27'
$ws.Range("V46").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/greuther-furth-hannover/zyA2QCuE/'

# Row 49
$ws.Range("F49").Value = 'Elversberg'
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 'Hamburger SV'
$ws.Range("I49").Value = 1
$ws.Range("J49").Value = 3.89
$ws.Range("K49").Value = '04/09/2023 08:42'
$ws.Range("L49").Value = 3.78
$ws.Range("M49").Value = '16/09/2023 12:59'
$ws.Range("N49").Value = 4.26
$ws.Range("O49").Value = '04/09/2023 08:42'
$ws.Range("P49").Value = 4.17
$ws.Range("Q49").Value = '16/09/2023 12:59'
$ws.Range("R49").Value = 1.85
$ws.Range("S49").Value = '04/09/2023 08:42'
$ws.Range("T49").Value = 1.9
$ws.Range("U49").Value = '16/09/2023 12:59'
$ws.Range("V49").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/elversberg-hamburger/6ux9XTns/'

# Row 50
$ws.Range("F50").Value = 'Hansa Rostock'
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 'Dusseldorf'
$ws.Range("I50").Value = 3
$ws.Range("J50").Value = 3.01
$ws.Range("K50").Value = '04/09/2023 08:42'
$ws.Range("L50").Value = 3.37
$ws.Range("M50").Value = '16/09/2023 12:59'
$ws.Range("N50").Value = 3.39
$ws.Range("O50").Value = '04/09/2023 08:42'
$ws.Range("P50").Value = 3.45
$ws.Range("Q50").Value = '16/09/2023 12:59'
$ws.Range("R50").Value = 2.43
$ws.Range("S50").Value = '04/09/2023 08:42'
$ws.Range("T50").Value = 2.26
$ws.Range("U50").Value = '16/09/2023 12:59'
$ws.Range("V50").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hansa-rostock-dusseldorf/zcowwjPJ/'

# Row 53
$ws.Range("F53").Value = 'St. Pauli'
$ws.Range("G53").Value = 5
$ws.Range("H53").Value = 'Holstein Kiel'
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = 1.65
$ws.Range("K53").Value = '04/09/2023 08:42'
$ws.Range("L53").Value = 1.93
$ws.Range("M53").Value = '17/09/2023 13:29'
$ws.Range("N53").Value = 4.28
$ws.Range("O53").Value = '04/09/2023 08:42'
$ws.Range("P53").Value = 3.7
$ws.Range("Q53").Value = '17/09/2023 13:29'
$ws.Range("R53").Value = 4.96
$ws.Range("S53").Value = '04/09/2023 08:42'
$ws.Range("T53").Value = 4.12
$ws.Range("U53").Value = '17/09/2023 13:28'
$ws.Range("V53").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/st-pauli-holstein-kiel/ny5rH820/'

# Row 54
$ws.Range("F54").Value = 'Hannover'
$ws.Range("G54").Value = 7
$ws.Range("H54").Value = 'VfL Osnabruck'
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 1.79
$ws.Range("K54").Value = '04/09/2023 08:42'
$ws.Range("L54").Value = 1.83
$ws.Range("M54").Value = '17/09/2023 13:28'
$ws.Range("N54").Value = 3.99
$ws.Range("O54").Value = '04/09/2023 08:42'
$ws.Range("P54").Value = 4.06
$ws.Range("Q54").Value = '17/09/2023 13:28'
$ws.Range("R54").Value = 4.3
$ws.Range("S54").Value = '04/09/2023 08:42'
$ws.Range("T54").Value = 4.19
$ws.Range("U54").Value = '17/09/2023 13:27'
$ws.Range("V54").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hannover-vfl-osnabruck/Yi7fE6nJ/'

# Row 55
$ws.Range("F55").Value = 'Hertha Berlin'
$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 'Braunschweig'
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1.59
$ws.Range("K55").Value = '04/09/2023 08:42'
$ws.Range("L55").Value = 1.78
$ws.Range("M55").Value = '17/09/2023 13:20'
$ws.Range("N55").Value = 4.48
$ws.Range("O55").Value = '04/09/2023 08:42'
$ws.Range("P55").Value = 4.2
$ws.Range("Q55").Value = '17/09/2023 13:29'
$ws.Range("R55").Value = 5.2
$ws.Range("S55").Value = '04/09/2023 08:42'
$ws.Range("T55").Value = 4.39
$ws.Range("U55").Value = '17/09/2023 13:23'
$ws.Range("V55").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hertha-berlin-braunschweig/6PmoIlmf/'

# Row 62
$ws.Range("F62").Value = 'Dusseldorf'
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 'Hannover'
$ws.Range("I62").Value = 1
$ws.Range("J62").Value = 1.71
$ws.Range("K62").Value = '17/09/2023 12:42'
$ws.Range("L62").Value = 1.98
$ws.Range("M62").Value = '24/09/2023 13:27'
$ws.Range("N62").Value = 4.24
$ws.Range("O62").Value = '17/09/2023 12:42'
$ws.Range("P62").Value = 3.87
$ws.Range("Q62").Value = '24/09/2023 13:28'
$ws.Range("R62").Value = 4.47
$ws.Range("S62").Value = '17/09/2023 12:42'
$ws.Range("T62").Value = 3.78
$ws.Range("U62").Value = '24/09/2023 13:27'
$ws.Range("V62").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/dusseldorf-hannover/pWzDWm2m/'

# Row 63
$ws.Range("F63").Value = 'Holstein Kiel'
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 'Hertha Berlin'
$ws.Range("I63").Value = 3
$ws.Range("J63").Value = 2.42
$ws.Range("K63").Value = '17/09/2023 16:13'
$ws.Range("L63").Value = 2.29
$ws.Range("M63").Value = '24/09/2023 13:20'
$ws.Range("N63").Value = 3.72
$ws.Range("O63").Value = '17/09/2023 16:13'
$ws.Range("P63").Value = 3.94
$ws.Range("Q63").Value = '24/09/2023 13:29'
$ws.Range("R63").Value = 2.8
$ws.Range("S63").Value = '17/09/2023 16:13'
$ws.Range("T63").Value = 2.97
$ws.Range("U63").Value = '24/09/2023 13:20'
$ws.Range("V63").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/holstein-kiel-hertha-berlin/vVTLURWa/'

# Row 64
$ws.Range("F64").Value = 'Kaiserslautern'
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = 'Hansa Rostock'
$ws.Range("I64").Value = 1
$ws.Range("J64").Value = 1.99
$ws.Range("K64").Value = '17/09/2023 12:42'
$ws.Range("L64").Value = 2.11
$ws.Range("M64").Value = '24/09/2023 13:22'
$ws.Range("N64").Value = 3.65
$ws.Range("O64").Value = '17/09/2023 12:42'
$ws.Range("P64").Value = 3.37
$ws.Range("Q64").Value = '24/09/2023 13:24'
$ws.Range("R64").Value = 3.9
$ws.Range("S64").Value = '17/09/2023 12:42'
$ws.Range("T64").Value = 3.88
$ws.Range("U64").Value = '24/09/2023 13:24'
$ws.Range("V64").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/kaiserslautern-hansa-rostock/YiTPTon6/'

# Row 65
$ws.Range("F65").Value = 'Hamburger SV'
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 'Dusseldorf'
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 1.74
$ws.Range("K65").Value = '24/09/2023 12:42'
$ws.Range("L65").Value = 1.97
$ws.Range("M65").Value = '29/09/2023 18:29'
$ws.Range("N65").Value = 4.51
$ws.Range("O65").Value = '24/09/2023 12:42'
$ws.Range("P65").Value = 4.13
$ws.Range("Q65").Value = '29/09/2023 18:29'
$ws.Range("R65").Value = 4.21
$ws.Range("S65").Value = '24/09/2023 12:42'
$ws.Range("T65").Value = 3.58
$ws.Range("U65").Value = '29/09/2023 18:29'
$ws.Range("V65").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hamburger-dusseldorf/KlkeqRog/'

# Row 66
$ws.Range("F66").Value = 'Paderborn'
$ws.Range("G66").Value = 3
$ws.Range("H66").Value = 'Schalke'
$ws.Range("I66").Value = 1
$ws.Range("J66").Value = 2.39
$ws.Range("K66").Value = '25/09/2023 11:42'
$ws.Range("L66").Value = 2.13
$ws.Range("M66").Value = '29/09/2023 18:29'
$ws.Range("N66").Value = 3.85
$ws.Range("O66").Value = '25/09/2023 11:42'
$ws.Range("P66").Value = 3.95
$ws.Range("Q66").Value = '29/09/2023 18:29'
$ws.Range("R66").Value = 2.78
$ws.Range("S66").Value = '25/09/2023 11:42'
$ws.Range("T66").Value = 3.27
$ws.Range("U66").Value = '29/09/2023 18:29'
$ws.Range("V66").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/paderborn-schalke/EP3bro1a/'

# Row 67
$ws.Range("F67").Value = 'Hannover'
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 'Wehen'
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1.69
$ws.Range("K67").Value = '25/09/2023 11:42'
$ws.Range("L67").Value = 1.78
$ws.Range("M67").Value = '30/09/2023 12:59'
$ws.Range("N67").Value = 4.2
$ws.Range("O67").Value = '25/09/2023 11:42'
$ws.Range("P67").Value = 4.17
$ws.Range("Q67").Value = '30/09/2023 12:59'
$ws.Range("R67").Value = 4.66
$ws.Range("S67").Value = '25/09/2023 11:42'
$ws.Range("T67").Value = 4.4
$ws.Range("U67").Value = '30/09/2023 12:59'
$ws.Range("V67").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hannover-wehen/zuh7tPVB/'

# Row 68
$ws.Range("F68").Value = 'Hansa Rostock'
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 'Braunschweig'
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2.1
$ws.Range("K68").Value = '24/09/2023 12:42'
$ws.Range("L68").Value = 2.18
$ws.Range("M68").Value = '30/09/2023 12:51'
$ws.Range("N68").Value = 3.58
$ws.Range("O68").Value = '24/09/2023 12:42'
$ws.Range("P68").Value = 3.56
$ws.Range("Q68").Value = '30/09/2023 12:57'
$ws.Range("R68").Value = 3.53
$ws.Range("S68").Value = '24/09/2023 12:42'
$ws.Range("T68").Value = 3.41
$ws.Range("U68").Value = '30/09/2023 12:38'
$ws.Range("V68").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hansa-rostock-braunschweig/lrdBuqoI/'

# Row 69
$ws.Range("F69").Value = 'Karlsruher SC'
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 'Holstein Kiel'
$ws.Range("I69").Value = 2
$ws.Range("J69").Value = 1.84
$ws.Range("K69").Value = '24/09/2023 12:42'
$ws.Range("L69").Value = 1.82
$ws.Range("M69").Value = '30/09/2023 12:59'
$ws.Range("N69").Value = 4.23
$ws.Range("O69").Value = '24/09/2023 12:42'
$ws.Range("P69").Value = 4.27
$ws.Range("Q69").Value = '30/09/2023 12:59'
$ws.Range("R69").Value = 3.83
$ws.Range("S69").Value = '24/09/2023 12:42'
$ws.Range("T69").Value = 4.05
$ws.Range("U69").Value = '30/09/2023 12:57'
$ws.Range("V69").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/karlsruher-holstein-kiel/hG22s5G5/'

# Row 71
$ws.Range("F71").Value = 'Elversberg'
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 'Greuther Furth'
$ws.Range("I71").Value = 1
$ws.Range("J71").Value = 2.7
$ws.Range("K71").Value = '25/09/2023 11:42'
$ws.Range("L71").Value = 2.48
$ws.Range("M71").Value = '01/10/2023 13:24'
$ws.Range("N71").Value = 3.76
$ws.Range("O71").Value = '25/09/2023 11:42'
$ws.Range("P71").Value = 3.64
$ws.Range("Q71").Value = '01/10/2023 13:27'
$ws.Range("R71").Value = 2.49
$ws.Range("S71").Value = '25/09/2023 11:42'
$ws.Range("T71").Value = 2.86
$ws.Range("U71").Value = '01/10/2023 13:16'
$ws.Range("V71").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/elversberg-greuther-furth/MefJwNFU/'

# Row 73
$ws.Range("F73").Value = 'VfL Osnabruck'
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 'Kaiserslautern'
$ws.Range("I73").Value = 2
$ws.Range("J73").Value = 2.76
$ws.Range("K73").Value = '25/09/2023 11:42'
$ws.Range("L73").Value = 2.62
$ws.Range("M73").Value = '01/10/2023 13:30'
$ws.Range("N73").Value = 3.54
$ws.Range("O73").Value = '25/09/2023 11:42'
$ws.Range("P73").Value = 3.59
$ws.Range("Q73").Value = '01/10/2023 13:04'
$ws.Range("R73").Value = 2.54
$ws.Range("S73").Value = '25/09/2023 11:42'
$ws.Range("T73").Value = 2.73
$ws.Range("U73").Value = '01/10/2023 13:30'
$ws.Range("V73").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/vfl-osnabruck-kaiserslautern/MVQPb6ot/'

# Row 80
$ws.Range("F80").Value = 'Holstein Kiel'
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 'Elversberg'
$ws.Range("I80").Value = 1
$ws.Range("J80").Value = 2.03
$ws.Range("K80").Value = '01/10/2023 12:43'
$ws.Range("L80").Value = 2.1
$ws.Range("M80").Value = '08/10/2023 13:28'
$ws.Range("N80").Value = 3.95
$ws.Range("O80").Value = '01/10/2023 12:43'
$ws.Range("P80").Value = 3.93
$ws.Range("Q80").Value = '08/10/2023 13:28'
$ws.Range("R80").Value = 3.39
$ws.Range("S80").Value = '01/10/2023 12:43'
$ws.Range("T80").Value = 3.36
$ws.Range("U80").Value = '08/10/2023 13:28'
$ws.Range("V80").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/holstein-kiel-elversberg/lrIteOp5/'

# Row 81
$ws.Range("F81").Value = 'Schalke'
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 'Hertha Berlin'
$ws.Range("I81").Value = 2
$ws.Range("J81").Value = 2.14
$ws.Range("K81").Value = '30/09/2023 19:43'
$ws.Range("L81").Value = 1.99
$ws.Range("M81").Value = '08/10/2023 13:29'
$ws.Range("N81").Value = 3.84
$ws.Range("O81").Value = '30/09/2023 19:43'
$ws.Range("P81").Value = 3.97
$ws.Range("Q81").Value = '08/10/2023 13:29'
$ws.Range("R81").Value = 3.3
$ws.Range("S81").Value = '30/09/2023 19:43'
$ws.Range("T81").Value = 3.63
$ws.Range("U81").Value = '08/10/2023 13:29'
$ws.Range("V81").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/schalke-hertha-berlin/0MPTcQ0n/'

# Row 82
$ws.Range("F82").Value = 'Braunschweig'
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 'Paderborn'
$ws.Range("I82").Value = 3
$ws.Range("J82").Value = 2.69
$ws.Range("K82").Value = '30/09/2023 12:43'
$ws.Range("L82").Value = 2.89
$ws.Range("M82").Value = '08/10/2023 13:27'
$ws.Range("N82").Value = 3.59
$ws.Range("O82").Value = '30/09/2023 12:43'
$ws.Range("P82").Value = 3.59
$ws.Range("Q82").Value = '08/10/2023 13:28'
$ws.Range("R82").Value = 2.58
$ws.Range("S82").Value = '30/09/2023 12:43'
$ws.Range("T82").Value = 2.48
$ws.Range("U82").Value = '08/10/2023 13:27'
$ws.Range("V82").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/braunschweig-paderborn/GdYMzfaj/'

# Row 89
$ws.Range("F89").Value = 'Nurnberg'
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 'Hertha Berlin'
$ws.Range("I89").Value = 1
$ws.Range("J89").Value = 2.65
$ws.Range("K89").Value = '10/10/2023 11:42'
$ws.Range("L89").Value = 2.56
$ws.Range("M89").Value = '22/10/2023 13:28'
$ws.Range("N89").Value = 3.7
$ws.Range("O89").Value = '10/10/2023 11:42'
$ws.Range("P89").Value = 3.55
$ws.Range("Q89").Value = '22/10/2023 13:28'
$ws.Range("R89").Value = 2.56
$ws.Range("S89").Value = '10/10/2023 11:42'
$ws.Range("T89").Value = 2.82
$ws.Range("U89").Value = '22/10/2023 13:28'
$ws.Range("V89").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/nurnberg-hertha-berlin/f3LlkgEq/'

# Row 90
$ws.Range("F90").Value = 'Hansa Rostock'
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 'Holstein Kiel'
$ws.Range("I90").Value = 3
$ws.Range("J90").Value = 2.46
$ws.Range("K90").Value = '09/10/2023 16:12'
$ws.Range("L90").Value = 2.65
$ws.Range("M90").Value = '22/10/2023 13:29'
$ws.Range("N90").Value = 3.66
$ws.Range("O90").Value = '09/10/2023 16:12'
$ws.Range("P90").Value = 3.29
$ws.Range("Q90").Value = '22/10/2023 13:29'
$ws.Range("R90").Value = 2.84
$ws.Range("S90").Value = '09/10/2023 16:12'
$ws.Range("T90").Value = 2.88
$ws.Range("U90").Value = '22/10/2023 13:29'
$ws.Range("V90").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hansa-rostock-holstein-kiel/lMrMgZ5S/'

# Row 91
$ws.Range("F91").Value = 'Karlsruher SC'
$ws.Range("G91").Value = 3
$ws.Range("H91").Value = 'Schalke'
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2.06
$ws.Range("K91").Value = '10/10/2023 11:42'
$ws.Range("L91").Value = 2.37
$ws.Range("M91").Value = '22/10/2023 13:29'
$ws.Range("N91").Value = 3.92
$ws.Range("O91").Value = '10/10/2023 11:42'
$ws.Range("P91").Value = 3.84
$ws.Range("Q91").Value = '22/10/2023 13:29'
$ws.Range("R91").Value = 3.41
$ws.Range("S91").Value = '10/10/2023 11:42'
$ws.Range("T91").Value = 2.9
$ws.Range("U91").Value = '22/10/2023 13:29'
$ws.Range("V91").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/karlsruher-schalke/tduEeeyG/'

# Row 94
$ws.Range("F94").Value = 'Hertha Berlin'
$ws.Range("G94").Value = 3
$ws.Range("H94").Value = 'Paderborn'
$ws.Range("I94").Value = 1
$ws.Range("J94").Value = 2.05
$ws.Range("K94").Value = '22/10/2023 13:42'
$ws.Range("L94").Value = 2.27
$ws.Range("M94").Value = '28/10/2023 13:00'
$ws.Range("N94").Value = 3.88
$ws.Range("O94").Value = '22/10/2023 13:42'
$ws.Range("P94").Value = 3.65
$ws.Range("Q94").Value = '28/10/2023 12:58'
$ws.Range("R94").Value = 3.48
$ws.Range("S94").Value = '22/10/2023 13:42'
$ws.Range("T94").Value = 3.19
$ws.Range("U94").Value = '28/10/2023 13:00'
$ws.Range("V94").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hertha-berlin-paderborn/veagmFb2/'

# Row 95
$ws.Range("F95").Value = 'Schalke'
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = 'Hannover'
$ws.Range("I95").Value = 2
$ws.Range("J95").Value = 1.83
$ws.Range("K95").Value = '22/10/2023 13:42'
$ws.Range("L95").Value = 2.39
$ws.Range("M95").Value = '28/10/2023 13:00'
$ws.Range("N95").Value = 4.19
$ws.Range("O95").Value = '22/10/2023 13:42'
$ws.Range("P95").Value = 3.77
$ws.Range("Q95").Value = '28/10/2023 13:00'
$ws.Range("R95").Value = 3.89
$ws.Range("S95").Value = '22/10/2023 13:42'
$ws.Range("T95").Value = 2.91
$ws.Range("U95").Value = '28/10/2023 13:00'
$ws.Range("V95").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/schalke-hannover/tKN0nib2/'

# Row 96
$ws.Range("F96").Value = 'St. Pauli'
$ws.Range("G96").Value = 2
$ws.Range("H96").Value = 'Karlsruher SC'
$ws.Range("I96").Value = 1
$ws.Range("J96").Value = 1.76
$ws.Range("K96").Value = '22/10/2023 13:42'
$ws.Range("L96").Value = 1.57
$ws.Range("M96").Value = '28/10/2023 12:58'
$ws.Range("N96").Value = 4.19
$ws.Range("O96").Value = '22/10/2023 13:42'
$ws.Range("P96").Value = 4.36
$ws.Range("Q96").Value = '28/10/2023 12:59'
$ws.Range("R96").Value = 4.39
$ws.Range("S96").Value = '22/10/2023 13:42'
$ws.Range("T96").Value = 6.01
$ws.Range("U96").Value = '28/10/2023 12:59'
$ws.Range("V96").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/st-pauli-karlsruher/xtIdnZD8/'

# Row 110
$ws.Range("F110").Value = 'St. Pauli'
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 'Hannover'
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 1.62
$ws.Range("K110").Value = '05/11/2023 13:42'
$ws.Range("L110").Value = 1.67
$ws.Range("M110").Value = '10/11/2023 18:07'
$ws.Range("N110").Value = 4.27
$ws.Range("O110").Value = '05/11/2023 13:42'
$ws.Range("P110").Value = 4.08
$ws.Range("Q110").Value = '10/11/2023 18:21'
$ws.Range("R110").Value = 5.19
$ws.Range("S110").Value = '05/11/2023 13:42'
$ws.Range("T110").Value = 5.34
$ws.Range("U110").Value = '10/11/2023 18:21'
$ws.Range("V110").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/st-pauli-hannover/bwDZmnkD/'

# Row 111
$ws.Range("F111").Value = 'Schalke'
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 'Elversberg'
$ws.Range("I111").Value = 2
$ws.Range("J111").Value = 1.96
$ws.Range("K111").Value = '04/11/2023 13:12'
$ws.Range("L111").Value = 2.05
$ws.Range("M111").Value = '10/11/2023 18:29'
$ws.Range("N111").Value = 4
$ws.Range("O111").Value = '04/11/2023 13:12'
$ws.Range("P111").Value = 3.96
$ws.Range("Q111").Value = '10/11/2023 18:29'
$ws.Range("R111").Value = 3.65
$ws.Range("S111").Value = '04/11/2023 13:12'
$ws.Range("T111").Value = 3.48
$ws.Range("U111").Value = '10/11/2023 18:29'
$ws.Range("V111").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/schalke-elversberg/htARk8K0/'

# Row 112
$ws.Range("F112").Value = 'Braunschweig'
$ws.Range("G112").Value = 3
$ws.Range("H112").Value = 'VfL Osnabruck'
$ws.Range("I112").Value = 2
$ws.Range("J112").Value = 2.16
$ws.Range("K112").Value = '05/11/2023 13:42'
$ws.Range("L112").Value = 2.45
$ws.Range("M112").Value = '11/11/2023 12:59'
$ws.Range("N112").Value = 3.79
$ws.Range("O112").Value = '05/11/2023 13:42'
$ws.Range("P112").Value = 3.38
$ws.Range("Q112").Value = '11/11/2023 12:59'
$ws.Range("R112").Value = 3.19
$ws.Range("S112").Value = '05/11/2023 13:42'
$ws.Range("T112").Value = 3.09
$ws.Range("U112").Value = '11/11/2023 12:59'
$ws.Range("V112").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/braunschweig-vfl-osnabruck/htWd5mKm/'

# Row 114
$ws.Range("F114").Value = 'Paderborn'
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 'Nurnberg'
$ws.Range("I114").Value = 3
$ws.Range("J114").Value = 1.95
$ws.Range("K114").Value = '05/11/2023 13:42'
$ws.Range("L114").Value = 2.27
$ws.Range("M114").Value = '11/11/2023 12:59'
$ws.Range("N114").Value = 3.96
$ws.Range("O114").Value = '05/11/2023 13:42'
$ws.Range("P114").Value = 3.86
$ws.Range("Q114").Value = '11/11/2023 12:58'
$ws.Range("R114").Value = 3.73
$ws.Range("S114").Value = '05/11/2023 13:42'
$ws.Range("T114").Value = 3.04
$ws.Range("U114").Value = '11/11/2023 12:59'
$ws.Range("V114").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/paderborn-nurnberg/0riICVRD/'

# Row 119
$ws.Range("F119").Value = 'Hannover'
$ws.Range("G119").Value = 2
$ws.Range("H119").Value = 'Hertha Berlin'
$ws.Range("I119").Value = 2
$ws.Range("J119").Value = 2.11
$ws.Range("K119").Value = '12/11/2023 15:42'
$ws.Range("L119").Value = 2.2
$ws.Range("M119").Value = '24/11/2023 18:08'
$ws.Range("N119").Value = 3.85
$ws.Range("O119").Value = '12/11/2023 15:42'
$ws.Range("P119").Value = 3.6
$ws.Range("Q119").Value = '24/11/2023 18:08'
$ws.Range("R119").Value = 3.26
$ws.Range("S119").Value = '12/11/2023 15:42'
$ws.Range("T119").Value = 3.39
$ws.Range("U119").Value = '24/11/2023 18:08'
$ws.Range("V119").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hannover-hertha-berlin/dnGG8MQu/'

# Row 120
$ws.Range("F120").Value = 'Hamburger SV'
$ws.Range("G120").Value = 2
$ws.Range("H120").Value = 'Braunschweig'
$ws.Range("I120").Value = 1
$ws.Range("J120").Value = 1.32
$ws.Range("K120").Value = '12/11/2023 15:42'
$ws.Range("L120").Value = 1.34
$ws.Range("M120").Value = '24/11/2023 18:10'
$ws.Range("N120").Value = 6.09
$ws.Range("O120").Value = '12/11/2023 15:42'
$ws.Range("P120").Value = 6.18
$ws.Range("Q120").Value = '24/11/2023 18:11'
$ws.Range("R120").Value = 8.01
$ws.Range("S120").Value = '12/11/2023 15:42'
$ws.Range("T120").Value = 8.09
$ws.Range("U120").Value = '24/11/2023 18:11'
$ws.Range("V120").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hamburger-braunschweig/06Z43Rl0/'

# Row 121
$ws.Range("F121").Value = 'Elversberg'
$ws.Range("G121").Value = 4
$ws.Range("H121").Value = 'Paderborn'
$ws.Range("I121").Value = 1
$ws.Range("J121").Value = 2.17
$ws.Range("K121").Value = '12/11/2023 15:42'
$ws.Range("L121").Value = 2.35
$ws.Range("M121").Value = '25/11/2023 12:58'
$ws.Range("N121").Value = 3.94
$ws.Range("O121").Value = '12/11/2023 15:42'
$ws.Range("P121").Value = 3.75
$ws.Range("Q121").Value = '25/11/2023 12:58'
$ws.Range("R121").Value = 3.14
$ws.Range("S121").Value = '12/11/2023 15:42'
$ws.Range("T121").Value = 2.98
$ws.Range("U121").Value = '25/11/2023 12:58'
$ws.Range("V121").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/elversberg-paderborn/W48T5KAb/'

# Row 122
$ws.Range("F122").Value = 'Hansa Rostock'
$ws.Range("G122").Value = 2
$ws.Range("H122").Value = 'St. Pauli'
$ws.Range("I122").Value = 3
$ws.Range("J122").Value = 4.5
$ws.Range("K122").Value = '12/11/2023 15:42'
$ws.Range("L122").Value = 5.28
$ws.Range("M122").Value = '25/11/2023 12:56'
$ws.Range("N122").Value = 4.02
$ws.Range("O122").Value = '12/11/2023 15:42'
$ws.Range("P122").Value = 3.83
$ws.Range("Q122").Value = '25/11/2023 12:59'
$ws.Range("R122").Value = 1.75
$ws.Range("S122").Value = '12/11/2023 15:42'
$ws.Range("T122").Value = 1.72
$ws.Range("U122").Value = '25/11/2023 12:56'
$ws.Range("V122").Value = 'https://www.betexplorer.com/football/germany/2-bundesliga/hansa-rostock-st-pauli/hE9P60eh/'

